$d = $word.ActiveDocument

# "CSV parsing" topic paragraph: tighten the wording that explains where the
# CSV (de)serialisation happens - drop the now-inaccurate "...and in the
# InnerRemoteSave/InnerRemoteLoad class of..." clause and join the two
# sentences with "and,".

$find1 = "values. This happens in the savePTSLocally method and in the InnerRemoteSave class of"
$repl1 = "values and, this happens in the savePTSLocally method of"
$found1 = $d.Content.Find.Execute($find1, $true, $false, $false, $false, $false, $true, 1, $false, $repl1, 2)

$find2 = "object. This happens in the localLoadPTS method and in the InnerRemoteLoad class of"
$repl2 = "object and, this happens in the localLoadPTS method of"
$found2 = $d.Content.Find.Execute($find2, $true, $false, $false, $false, $false, $true, 1, $false, $repl2, 2)

Write-Output ("replace1=" + $found1 + " replace2=" + $found2)
